$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5 (pushes the existing row 5 and everything below it down by one).
$ws.Rows("5:5").Insert()

# The row-insert operation copies formatting down from row 4 into the blank new row 5 for
# columns F:O and U:V (which had styled-but-empty cells in row 4). The new "MasterCard"
# record row does not carry any of that inherited formatting, so clear it back out.
$ws.Range("F5:O5").Clear()
$ws.Range("U5:V5").Clear()

# Populate the new payment-details row, in the same left-to-right-ish entry order used by the
# author (cardType, cardNumber, then the row "name" cell) so new shared-string entries land in
# the same order as the target workbook.
$ws.Range("W5").Value = "MasterCard"

$cardNumber = $ws.Range("X5")
$cardNumber.Value = "'5111005111051128"
$cardNumber.NumberFormat = "@"
$cardNumber.HorizontalAlignment = -4131

$ws.Range("Y5").Value = 2025
$ws.Range("Y5").NumberFormat = "0"

$ws.Range("Z5").Value = "Feb"

$ws.Range("AB5").Value = 123
$ws.Range("AB5").NumberFormat = "0;[Red]0"

$ws.Range("A5").Value = "Ccmastercard"

# Column W now holds "MasterCard" which needs a widened, best-fit column.
$ws.Columns("W:W").AutoFit()

# Move the selection/cursor to the newly added row, and let the view scroll back to the top.
[void]$ws.Range("A5").Select()
